$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ratio_threshold_range row (row 4): Min 0.8 -> 0.7, Max 1.4 -> 1.3
$ws.Range("B4").Value = 0.7
$ws.Range("C4").Value = 1.3

# The theta_threshold_range row (row 5) is being replaced by the
# pie_threshold_range row's label/values, and the old row 6
# (pie_threshold_range) is removed - net effect: delete row 5
# (theta_threshold_range) so that the old row 6 shifts up to become row 5.
$ws.Rows("5").Delete()

# Update the (now) row 5 values (previously row 6, pie_threshold_range)
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Update selection to match target state
$ws.Range("E6").Select()

# Match the page setup that was applied (paper size / orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
